# Updated cryptos list - apply latest price/volume figures scraped from
# coinranking.com. Rows 44/45 and 48/49 also swap their coin identity
# (name/link/price/volume) to reflect the new ranking order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.730.60"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +4.32%  "

# Row 3 - Ethereum
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.074.89"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.82%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.09%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.89%  "

# Row 6 - Solana
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "141.66"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.41%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  -0.05%  "

# Row 8 - LidoStakedEther
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.063.57"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +2.87%  "

# Row 9 - XRP
$ws.Range("E9").Value = "  +1.23%  "

# Row 10 - Dogecoin
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.139"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.67%  "

# Row 11 - Toncoin
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.55"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +9.16%  "

# Row 12 - Cardano
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.466"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.56%  "

# Row 13 - ShibaInu
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000238"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.70%  "

# Row 14 - Avalanche
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.23"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +4.73%  "

# Row 15 - TRON
$ws.Range("E15").Value = "  +0.19%  "

# Row 16 - WrappedliquidstakedEther2.0
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.580.54"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.79%  "

# Row 17 - Polkadot
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "7.25"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +3.00%  "

# Row 18 - WrappedEther
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.070.21"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.86%  "

# Row 19 - WrappedBTC
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "61.665.23"
$ws.Range("D19").Style = "Normal"

# Row 20 - BitcoinCash
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "447.93"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +4.88%  "

# Row 21 - Chainlink
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.89"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.69%  "

# Row 22 - Polygon
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.728"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.00%  "

# Row 23 - Uniswap
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.42"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.86%  "

# Row 24 - InternetComputer(DFINITY)
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.76"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.63%  "

# Row 25 - Litecoin
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "81.86"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.61%  "

# Row 26 - Dai
$ws.Range("E26").Value = "  -0.08%  "

# Row 27 - ImmutableX
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.26"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.25%  "

# Row 28 - FirstDigitalUSD
$ws.Range("E28").Value = "  +0.20%  "

# Row 29 - PancakeSwap
$ws.Range("E29").Value = "  +4.70%  "

# Row 30 - RenderToken
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.10"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.86%  "

# Row 31 - NEARProtocol
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.72"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +9.99%  "

# Row 32 - Hedera
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.114"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +15.34%  "

# Row 33 - EthereumClassic
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "26.65"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +3.85%  "

# Row 34 - Mantle
$ws.Range("E34").Value = "  +4.83%  "

# Row 35 - PEPE
$ws.Range("E35").Value = "  +2.61%  "

# Row 36 - Filecoin
$ws.Range("E36").Value = "  +3.36%  "

# Row 37 - Stacks
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.18"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +5.66%  "

# Row 38 - OKB
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "50.40"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.54%  "

# Row 39 - dogwifhat
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.99"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +8.56%  "

# Row 40 - Cosmos
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "8.82"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.19%  "

# Row 41 - Bittensor
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "427.62"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +7.03%  "

# Row 42 - VeChain
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0371"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +5.63%  "

# Row 43 - Maker
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.795.75"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.65%  "

# Row 44 - was Kaspa, now TheGraph (rows 44/45 swapped)
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.270"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +7.12%  "

# Row 45 - was TheGraph, now Kaspa
$ws.Range("B45").Value = "Kaspa"
$ws.Range("C45").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.109"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.31%  "

# Row 46 - Fetch.AI
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.10"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +5.11%  "

# Row 47 - USDe
$ws.Range("E47").Value = "  -0.03%  "

# Row 48 - was Monero, now Arweave (rows 48/49 swapped)
$ws.Range("B48").Value = "Arweave"
$ws.Range("C48").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "35.10"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +6.78%  "

# Row 49 - was Arweave, now Monero
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "123.85"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.41%  "

# Row 50 - Stellar
$ws.Range("E50").Value = "  +1.16%  "

# Row 51 - InjectiveProtocol
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "24.02"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.44%  "

